# "generacja recepty w formacie pdf"
#
# Logs a new batch of work done on 2025-08-05 (serial 45874) into the
# first empty block of rows (38-44) of the time-tracking table in
# columns B (date) / C (file) / D (lines changed), and moves the
# viewport/selection to where the new rows were entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 39-44 don't have any pre-existing formatting on column B (only
# B38 does, carried over from the template), so copy B38's date format
# down the column before writing the new dates into it.
$ws.Range("B38").Copy() | Out-Null
$ws.Range("B39:B44").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# File names were typed in this order (matches the order new entries
# were added to the workbook), each with its associated date + line count.
$entries = @(
    @{ Row = 38; File = "PdfGenratorService";       Lines = 28 },
    @{ Row = 44; File = "wykonaneBadaniaService";   Lines = 4 },
    @{ Row = 41; File = "wykonaneBadania.cs";       Lines = 2 },
    @{ Row = 42; File = "wykonaneBadaniaDTO";       Lines = 2 },
    @{ Row = 43; File = "DBInit.cs";                Lines = 4 },
    @{ Row = 39; File = "wykonaneBadaniaControler"; Lines = 25 },
    @{ Row = 40; File = "Program.cs";               Lines = 2 }
)

foreach ($entry in $entries) {
    $r = $entry.Row
    $ws.Range("B$r").Value = 45874
    $ws.Range("C$r").Value = $entry.File
    $ws.Range("D$r").Value = $entry.Lines
}

# Move the visible window / selection down to the freshly entered rows.
$ws.Range("F40").Select()
